$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.9105276959072057
$arr[0,1] = 0.1289059752690278
$arr[0,2] = 0
$arr[0,3] = 0.1681665655645261
$arr[0,4] = 2.908886162381549
$arr[0,5] = 1.833931374225813
$arr[0,6] = 1.591573240477217
$arr[0,7] = 0
$arr[0,8] = 0.1243362763011104
$arr[0,9] = 0.4391681502629865
$arr[0,10] = 0.3355767958224476
$ws.Range("B2:L2").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.876344550185479
$arr[0,1] = 0.1280678409840661
$arr[0,2] = 0
$arr[0,3] = 0.1675709846273108
$arr[0,4] = 2.898344445644398
$arr[0,5] = 1.829596865680884
$arr[0,6] = 1.594610352395051
$arr[0,7] = 0
$arr[0,8] = 0.1247942283774099
$arr[0,9] = 0.4090509805716351
$arr[0,10] = 0.3298141176898355
$ws.Range("B3:L3").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.8558412150274251
$arr[0,1] = 0.1275401981906725
$arr[0,2] = 0
$arr[0,3] = 0.1672672887482243
$arr[0,4] = 2.893139092647516
$arr[0,5] = 1.827767304001313
$arr[0,6] = 1.597041861104117
$arr[0,7] = 0
$arr[0,8] = 0.1250854085177577
$arr[0,9] = 0.3907772885610967
$arr[0,10] = 0.3264355330072135
$ws.Range("B4:L4").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.8476083872305367
$arr[0,1] = 0.1273218951867925
$arr[0,2] = 0
$arr[0,3] = 0.1671591665797223
$arr[0,4] = 2.891336627530237
$arr[0,5] = 1.827230699284783
$arr[0,6] = 1.598175239731461
$arr[0,7] = 0
$arr[0,8] = 0.1252065843288523
$arr[0,9] = 0.3833857529157285
$arr[0,10] = 0.3250990359969848
$ws.Range("B5:L5").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.8462487406710295
$arr[0,1] = 0.1272854475149785
$arr[0,2] = 0
$arr[0,3] = 0.1671421587809014
$arr[0,4] = 2.891056584358907
$arr[0,5] = 1.827154210597598
$arr[0,6] = 1.598372044664529
$arr[0,7] = 0
$arr[0,8] = 0.1252268576535389
$arr[0,9] = 0.3821617314089849
$arr[0,10] = 0.3248795500738595
$ws.Range("B6:L6").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.8557296878763054
$arr[0,1] = 0.1275372673894957
$arr[0,2] = 0
$arr[0,3] = 0.1672657672000071
$arr[0,4] = 2.893113493225684
$arr[0,5] = 1.827759221400314
$arr[0,6] = 1.597056569175976
$arr[0,7] = 0
$arr[0,8] = 0.1250870325383779
$arr[0,9] = 0.3906773801483752
$arr[0,10] = 0.3264173451656376
$ws.Range("B7:L7").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.8986408741031653
$arr[0,1] = 0.128619684244331
$arr[0,2] = 0
$arr[0,3] = 0.167948370566851
$arr[0,4] = 2.904988390327063
$arr[0,5] = 1.832264050410913
$arr[0,6] = 1.592502823023423
$arr[0,7] = 0
$arr[0,8] = 0.1244921060335256
$arr[0,9] = 0.4287385540537514
$arr[0,10] = 0.333556754969095
$ws.Range("B8:L8").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.9866267246038376
$arr[0,1] = 0.1306394333851912
$arr[0,2] = 0
$arr[0,3] = 0.1697769743733417
$arr[0,4] = 2.938331264062384
$arr[0,5] = 1.847710540038662
$arr[0,6] = 1.588069907156793
$arr[0,7] = 0
$arr[0,8] = 0.1234045642069277
$arr[0,9] = 0.5051033847910276
$arr[0,10] = 0.3488194228217765
$ws.Range("B9:L9").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.053599795163592
$arr[0,1] = 0.1320614003870162
$arr[0,2] = 0
$arr[0,3] = 0.1714170028822757
$arr[0,4] = 2.96896647282297
$arr[0,5] = 1.863110667007192
$arr[0,6] = 1.587556593163072
$arr[0,7] = 0
$arr[0,8] = 0.1226534586224712
$arr[0,9] = 0.5622603510783506
$arr[0,10] = 0.3607973564195674
$ws.Range("B10:L10").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.084572105687954
$arr[0,1] = 0.1326950122750716
$arr[0,2] = 0
$arr[0,3] = 0.1722270523354581
$arr[0,4] = 2.984238168224877
$arr[0,5] = 1.871000965257309
$arr[0,6] = 1.587919380379077
$arr[0,7] = 0
$arr[0,8] = 0.1223220980913728
$arr[0,9] = 0.5884911498303893
$arr[0,10] = 0.3664114666553644
$ws.Range("B11:L11").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.096372979181382
$arr[0,1] = 0.1329330512386449
$arr[0,2] = 0
$arr[0,3] = 0.1725429567639623
$arr[0,4] = 2.990213260911219
$arr[0,5] = 1.874116341184759
$arr[0,6] = 1.588142528265934
$arr[0,7] = 0
$arr[0,8] = 0.1221981006454111
$arr[0,9] = 0.5984569984952941
$arr[0,10] = 0.3685610397621275
$ws.Range("B12:L12").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.093828239588191
$arr[0,1] = 0.1328818694757388
$arr[0,2] = 0
$arr[0,3] = 0.1724745145901849
$arr[0,4] = 2.988917877419084
$arr[0,5] = 1.873439714700964
$arr[0,6] = 1.588090654457801
$arr[0,7] = 0
$arr[0,8] = 0.1222247398591874
$arr[0,9] = 0.5963092188891892
$arr[0,10] = 0.3680970418866991
$ws.Range("B13:L13").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.085541523404004
$arr[0,1] = 0.1327146338319452
$arr[0,2] = 0
$arr[0,3] = 0.1722528587382577
$arr[0,4] = 2.9847258941939
$arr[0,5] = 1.871254712439082
$arr[0,6] = 1.587936019841038
$arr[0,7] = 0
$arr[0,8] = 0.1223118670679648
$arr[0,9] = 0.5893103898606853
$arr[0,10] = 0.366587840499534
$ws.Range("B14:L14").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.080475077817937
$arr[0,1] = 0.1326119505113965
$arr[0,2] = 0
$arr[0,3] = 0.1721182791947165
$arr[0,4] = 2.982183191146405
$arr[0,5] = 1.869932946651033
$arr[0,6] = 1.587852471881035
$arr[0,7] = 0
$arr[0,8] = 0.1223654278722286
$arr[0,9] = 0.5850276715947302
$arr[0,10] = 0.3656664845370443
$ws.Range("B15:L15").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.051585826340613
$arr[0,1] = 0.1320197265370062
$arr[0,2] = 0
$arr[0,3] = 0.1713653480180497
$arr[0,4] = 2.96799530506803
$arr[0,5] = 1.862612843204857
$arr[0,6] = 1.58754488499801
$arr[0,7] = 0
$arr[0,8] = 0.1226753212772516
$arr[0,9] = 0.5605507111266377
$arr[0,10] = 0.3604337754104279
$ws.Range("B16:L16").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.033992485553512
$arr[0,1] = 0.1316530294939433
$arr[0,2] = 0
$arr[0,3] = 0.1709198064535222
$arr[0,4] = 2.959633555764086
$arr[0,5] = 1.858348984163399
$arr[0,6] = 1.587508931331257
$arr[0,7] = 0
$arr[0,8] = 0.1228680724650291
$arr[0,9] = 0.5455935820094169
$arr[0,10] = 0.35726591439402
$ws.Range("B17:L17").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.023920906067673
$arr[0,1] = 0.1314408675447716
$arr[0,2] = 0
$arr[0,3] = 0.1706695687010864
$arr[0,4] = 2.954949811613858
$arr[0,5] = 1.855979781412145
$arr[0,6] = 1.587544376610992
$arr[0,7] = 0
$arr[0,8] = 0.1229799094708417
$arr[0,9] = 0.537012281177681
$arr[0,10] = 0.3554594103410977
$ws.Range("B18:L18").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.020519038219874
$arr[0,1] = 0.1313688186222635
$arr[0,2] = 0
$arr[0,3] = 0.1705858788617292
$arr[0,4] = 2.953385568336898
$arr[0,5] = 1.855191901118246
$arr[0,6] = 1.58756601663859
$arr[0,7] = 0
$arr[0,8] = 0.1230179425177491
$arr[0,9] = 0.5341105247030384
$arr[0,10] = 0.3548504367395964
$ws.Range("B19:L19").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.035860398712316
$arr[0,1] = 0.1316921940691529
$arr[0,2] = 0
$arr[0,3] = 0.1709666117808517
$arr[0,4] = 2.96051066759793
$arr[0,5] = 1.858794260556721
$arr[0,6] = 1.58750694978562
$arr[0,7] = 0
$arr[0,8] = 0.122847453231592
$arr[0,9] = 0.5471835555325129
$arr[0,10] = 0.3576015285991758
$ws.Range("B20:L20").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.087973573079296
$arr[0,1] = 0.132763806382755
$arr[0,2] = 0
$arr[0,3] = 0.1723177162992258
$arr[0,4] = 2.985951969735169
$arr[0,5] = 1.871893038326704
$arr[0,6] = 1.587979111875654
$arr[0,7] = 0
$arr[0,8] = 0.1222862355044714
$arr[0,9] = 0.5913652268527585
$arr[0,10] = 0.3670304892659146
$ws.Range("B21:L21").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.122453958310416
$arr[0,1] = 0.1334531192588244
$arr[0,2] = 0
$arr[0,3] = 0.1732540900018513
$arr[0,4] = 3.003698637213247
$arr[0,5] = 1.881197089321148
$arr[0,6] = 1.588787630563758
$arr[0,7] = 0
$arr[0,8] = 0.1219280820547972
$arr[0,9] = 0.6204316661130918
$arr[0,10] = 0.3733305393323718
$ws.Range("B22:L22").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.104012705816274
$arr[0,1] = 0.1330862278815061
$arr[0,2] = 0
$arr[0,3] = 0.1727494633499411
$arr[0,4] = 2.994124497480854
$arr[0,5] = 1.876163245657835
$arr[0,6] = 1.588310357680626
$arr[0,7] = 0
$arr[0,8] = 0.1221184460235953
$arr[0,9] = 0.6049009445721367
$arr[0,10] = 0.3699555313061182
$ws.Range("B23:L23").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 1.035015780994655
$arr[0,1] = 0.1316744919497808
$arr[0,2] = 0
$arr[0,3] = 0.1709454326804227
$arr[0,4] = 2.96011374051696
$arr[0,5] = 1.858592695229078
$arr[0,6] = 1.587507670836374
$arr[0,7] = 0
$arr[0,8] = 0.122856772003388
$arr[0,9] = 0.5464646731916787
$arr[0,10] = 0.3574497514909041
$ws.Range("B24:L24").Value = $arr

$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 0.9624144152952283
$arr[0,1] = 0.1301039714613914
$arr[0,2] = 0
$arr[0,3] = 0.1692300545969907
$arr[0,4] = 2.928233756356803
$arr[0,5] = 1.842821827432658
$arr[0,6] = 1.588787485082506
$arr[0,7] = 0
$arr[0,8] = 0.1236903371600371
$arr[0,9] = 0.4842599168963773
$arr[0,10] = 0.3445558581197616
$ws.Range("B25:L25").Value = $arr
